$wb = $excel.ActiveWorkbook

# --- Sheet 1 "full list": set B11 (gate drive row) to FAN7888 ---
$ws1 = $wb.Worksheets.Item("full list")
$ws1.Range("B11").Value = "FAN7888"
$ws1.Range("B12").Select() | Out-Null

# --- Sheet 2 "top_pick": insert a new row at 4 for "gate drive A" / FAN7888,
#     relabel the old "gate drive" row (now row 5) as "gate drive B" ---
$ws2 = $wb.Worksheets.Item("top_pick")
$ws2.Rows.Item(4).Insert()

$ws2.Range("A5").Value = "gate drive B"

$ws2.Range("A4").Value = "gate drive A"
$ws2.Range("B4").Value = "FAN7888"
$ws2.Range("C4").Value = "custom"
$ws2.Range("D4").Value = "custom"
$ws2.Range("E4").Value = "yes"

$ws2.Range("H4").Select() | Out-Null

# --- New sheet "F28027 pin assignment" with PIN header and values 1-26 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "F28027 pin assignment"

$ws3.Range("A1").Value = "PIN"
for ($i = 1; $i -le 26; $i++) {
    $ws3.Cells.Item($i + 1, 1).Value = $i
}

$ws3.Range("A2:A27").Select() | Out-Null
